# Updates the cryptos worksheet with refreshed price/volume figures,
# and swaps the WrappedEther/Polkadot row order (rows 18-19) to match
# the latest GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.452.99'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.945.88'
$ws.Range("E3").Value = '  +2.16%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.16'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.01'
$ws.Range("E6").Value = '  +4.57%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.946.43'
$ws.Range("E8").Value = '  +2.19%  '
$ws.Range("E9").Value = '  +3.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.94'
$ws.Range("E11").Value = '  +9.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.434'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("E13").Value = '  +6.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.14'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.434.45'
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.455.30'
$ws.Range("E17").Value = '  +2.37%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.65'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.946.16'
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '433.91'
$ws.Range("E20").Value = '  +2.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.40'
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.660'
$ws.Range("E22").Value = '  +1.45%  '
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.14'
$ws.Range("E24").Value = '  +6.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.05'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.86'
$ws.Range("E26").Value = '  +4.91%  '
$ws.Range("E27").Value = '  +2.34%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.17'
$ws.Range("E29").Value = '  +7.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.16'
$ws.Range("E30").Value = '  +4.16%  '
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("E32").Value = '  +17.83%  '
$ws.Range("E33").Value = '  +3.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.10'
$ws.Range("E34").Value = '  +1.68%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  +2.27%  '
$ws.Range("E37").Value = '  +2.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.00'
$ws.Range("E38").Value = '  +7.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.57'
$ws.Range("E39").Value = '  +1.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("E40").Value = '  +5.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.33'
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.114'
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.273'
$ws.Range("E43").Value = '  +3.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.17'
$ws.Range("E44").Value = '  -1.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '134.88'
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.672.71'
$ws.Range("E46").Value = '  +0.86%  '
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '353.67'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.50'
$ws.Range("E51").Value = '  -0.11%  '
